# akaunting 3.0 (the last dance)
# categories.xlsx -- "enabled" column (D) was imported as native Excel
# booleans (TRUE/FALSE). Re-author it as literal text "TRUE" (so it
# round-trips through the CSV/XLSX importer as a string) and add a
# TRUE/FALSE list data-validation to the column so future edits stay
# constrained to those two values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-type column D (rows 2-7, "enabled") as the literal text "TRUE" ---
#
# A plain  $ws.Range("D2").Value = "TRUE"  gets auto-recognised as the
# boolean literal TRUE by Excel's input parser, same as before. Entering
# it with a leading apostrophe forces text, but leaves a "quote prefix"
# style behind on the cell. Instead, stage the forced-text value in a
# scratch cell, copy it, and paste only the *value* (PasteSpecial
# xlPasteValues = -4163) into each target cell -- this carries over the
# text type without carrying over the quote-prefix formatting.
$scratch = $ws.Range("F8")
$scratch.Value = "'TRUE"
$scratch.Copy()

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 4).PasteSpecial(-4163)
}

# Clean up the scratch cell, leaving the cursor where it was left.
$scratch.Clear()
$scratch.Select()

# --- 2. Constrain column D to a TRUE/FALSE dropdown ---
$rng = $ws.Range("D2:D1048576")
$rng.Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$rng.Validation.ErrorTitle = "Enabled Error"
$rng.Validation.ErrorMessage = "You must choose true or false"
